# Commit: "commit 2 par ghada"
#
# The paragraph currently reads "Salut 3A11". The author selected the
# word "Salut" and retyped it as "Bonjour" (Word drops a collapsed
# "_GoBack" bookmark right after the freshly-typed text to remember the
# last edit point), and then made a second small edit appending a
# trailing space at the very end of the document. Both edits leave the
# text split across three runs with the bookmark sitting between the
# first and second run.

$d = $word.ActiveDocument

# --- Edit 1: replace "Salut" with "Bonjour" -------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
# Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,
#         MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)
$found = $rng.Find.Execute("Salut", $true, $true, $false, $false, $false, `
                            $true, 1, $false, "", 0)

if ($found) {
    # Typing over the found selection replaces its text...
    $rng.Text = "Bonjour"

    # ...and Word leaves the collapsed "_GoBack" bookmark right after the
    # text that was just typed, marking the last edit location.
    $goBackPoint = $d.Range($rng.End, $rng.End)
    $d.Bookmarks.Add("_GoBack", $goBackPoint)
}

# --- Edit 2: a trailing space typed at the very end of the document ------
$docEnd = $d.Content.End
$endPoint = $d.Range($docEnd - 1, $docEnd - 1)
$endPoint.InsertAfter(" ")
